$wb = $excel.ActiveWorkbook

# Year 2000 (sheet index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 3).Value = 5943366
$ws.Cells.Item(2, 4).Value = 205.73316955566406
$ws.Cells.Item(2, 5).Value = 0.97245633602142334
$ws.Cells.Item(3, 3).Value = 5943366
$ws.Cells.Item(3, 4).Value = 205.73316955566406
$ws.Cells.Item(3, 5).Value = 0.97516489028930664
$ws.Cells.Item(4, 3).Value = 5943366
$ws.Cells.Item(4, 4).Value = 205.73316955566406
$ws.Cells.Item(4, 5).Value = 0.98107337951660156
$ws.Cells.Item(5, 3).Value = 5943366
$ws.Cells.Item(5, 4).Value = 205.73316955566406
$ws.Cells.Item(5, 5).Value = 0.99757492542266846
$ws.Cells.Item(6, 3).Value = 5943366
$ws.Cells.Item(6, 4).Value = 205.73316955566406
$ws.Cells.Item(6, 5).Value = 0.99952518939971924
$ws.Cells.Item(7, 3).Value = 5943366
$ws.Cells.Item(7, 4).Value = 205.73316955566406
$ws.Cells.Item(7, 5).Value = 0.99990880489349365
$ws.Cells.Item(8, 3).Value = 5943366
$ws.Cells.Item(8, 4).Value = 205.73316955566406
$ws.Cells.Item(8, 5).Value = 0.99998670816421509
$ws.Cells.Item(9, 3).Value = 5943366
$ws.Cells.Item(9, 4).Value = 205.73316955566406
$ws.Cells.Item(9, 5).Value = 0.99999243021011353
$ws.Cells.Item(10, 3).Value = 5943366
$ws.Cells.Item(10, 4).Value = 205.73316955566406
$ws.Cells.Item(10, 5).Value = 0.99999964237213135
$ws.Cells.Item(11, 3).Value = 5943366
$ws.Cells.Item(11, 4).Value = 205.73316955566406

# Year 2001 (sheet index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 3).Value = 5968060
$ws.Cells.Item(2, 4).Value = 209.95283508300781
$ws.Cells.Item(2, 5).Value = 0.97199761867523193
$ws.Cells.Item(3, 3).Value = 5968060
$ws.Cells.Item(3, 4).Value = 209.95283508300781
$ws.Cells.Item(3, 5).Value = 0.97491896152496338
$ws.Cells.Item(4, 3).Value = 5968060
$ws.Cells.Item(4, 4).Value = 209.95283508300781
$ws.Cells.Item(4, 5).Value = 0.9808925986289978
$ws.Cells.Item(5, 3).Value = 5968060
$ws.Cells.Item(5, 4).Value = 209.95283508300781
$ws.Cells.Item(5, 5).Value = 0.99738472700119019
$ws.Cells.Item(6, 3).Value = 5968060
$ws.Cells.Item(6, 4).Value = 209.95283508300781
$ws.Cells.Item(6, 5).Value = 0.99949163198471069
$ws.Cells.Item(7, 3).Value = 5968060
$ws.Cells.Item(7, 4).Value = 209.95283508300781
$ws.Cells.Item(7, 5).Value = 0.99991267919540405
$ws.Cells.Item(8, 3).Value = 5968060
$ws.Cells.Item(8, 4).Value = 209.95283508300781
$ws.Cells.Item(8, 5).Value = 0.9999842643737793
$ws.Cells.Item(9, 3).Value = 5968060
$ws.Cells.Item(9, 4).Value = 209.95283508300781
$ws.Cells.Item(9, 5).Value = 0.999991774559021
$ws.Cells.Item(10, 3).Value = 5968060
$ws.Cells.Item(10, 4).Value = 209.95283508300781
$ws.Cells.Item(11, 3).Value = 5968060
$ws.Cells.Item(11, 4).Value = 209.95283508300781

# Year 2002 (sheet index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 3).Value = 5986631
$ws.Cells.Item(2, 4).Value = 204.59521484375
$ws.Cells.Item(2, 5).Value = 0.9721110463142395
$ws.Cells.Item(3, 3).Value = 5986631
$ws.Cells.Item(3, 4).Value = 204.59521484375
$ws.Cells.Item(3, 5).Value = 0.97517186403274536
$ws.Cells.Item(4, 3).Value = 5986631
$ws.Cells.Item(4, 4).Value = 204.59521484375
$ws.Cells.Item(4, 5).Value = 0.98161154985427856
$ws.Cells.Item(5, 3).Value = 5986631
$ws.Cells.Item(5, 4).Value = 204.59521484375
$ws.Cells.Item(5, 5).Value = 0.99745702743530273
$ws.Cells.Item(6, 3).Value = 5986631
$ws.Cells.Item(6, 4).Value = 204.59521484375
$ws.Cells.Item(6, 5).Value = 0.99950557947158813
$ws.Cells.Item(7, 3).Value = 5986631
$ws.Cells.Item(7, 4).Value = 204.59521484375
$ws.Cells.Item(7, 5).Value = 0.99991196393966675
$ws.Cells.Item(8, 3).Value = 5986631
$ws.Cells.Item(8, 4).Value = 204.59521484375
$ws.Cells.Item(8, 5).Value = 0.99998563528060913
$ws.Cells.Item(9, 3).Value = 5986631
$ws.Cells.Item(9, 4).Value = 204.59521484375
$ws.Cells.Item(9, 5).Value = 0.99999183416366577
$ws.Cells.Item(10, 3).Value = 5986631
$ws.Cells.Item(10, 4).Value = 204.59521484375
$ws.Cells.Item(10, 5).Value = 0.99999964237213135

# Year 2003 (sheet index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 3).Value = 5998599
$ws.Cells.Item(2, 4).Value = 218.30923461914062
$ws.Cells.Item(2, 5).Value = 0.97003334760665894
$ws.Cells.Item(3, 3).Value = 5998599
$ws.Cells.Item(3, 4).Value = 218.30923461914062
$ws.Cells.Item(3, 5).Value = 0.97344982624053955
$ws.Cells.Item(4, 3).Value = 5998599
$ws.Cells.Item(4, 4).Value = 218.30923461914062
$ws.Cells.Item(4, 5).Value = 0.98015570640563965
$ws.Cells.Item(5, 3).Value = 5998599
$ws.Cells.Item(5, 4).Value = 218.30923461914062
$ws.Cells.Item(5, 5).Value = 0.99733239412307739
$ws.Cells.Item(6, 3).Value = 5998599
$ws.Cells.Item(6, 4).Value = 218.30923461914062
$ws.Cells.Item(6, 5).Value = 0.99946653842926025
$ws.Cells.Item(7, 3).Value = 5998599
$ws.Cells.Item(7, 4).Value = 218.30923461914062
$ws.Cells.Item(7, 5).Value = 0.99990630149841309
$ws.Cells.Item(8, 3).Value = 5998599
$ws.Cells.Item(8, 4).Value = 218.30923461914062
$ws.Cells.Item(8, 5).Value = 0.99998533725738525
$ws.Cells.Item(9, 3).Value = 5998599
$ws.Cells.Item(9, 4).Value = 218.30923461914062
$ws.Cells.Item(9, 5).Value = 0.99999183416366577
$ws.Cells.Item(10, 3).Value = 5998599
$ws.Cells.Item(10, 4).Value = 218.30923461914062
$ws.Cells.Item(10, 5).Value = 0.99999964237213135

# Year 2004 (sheet index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 3).Value = 6004671
$ws.Cells.Item(2, 4).Value = 230.83682250976562
$ws.Cells.Item(2, 5).Value = 0.96806567907333374
$ws.Cells.Item(3, 3).Value = 6004671
$ws.Cells.Item(3, 4).Value = 230.83682250976562
$ws.Cells.Item(3, 5).Value = 0.97190040349960327
$ws.Cells.Item(4, 3).Value = 6004671
$ws.Cells.Item(4, 4).Value = 230.83682250976562
$ws.Cells.Item(4, 5).Value = 0.97897803783416748
$ws.Cells.Item(5, 3).Value = 6004671
$ws.Cells.Item(5, 4).Value = 230.83682250976562
$ws.Cells.Item(5, 5).Value = 0.99717605113983154
$ws.Cells.Item(6, 3).Value = 6004671
$ws.Cells.Item(6, 4).Value = 230.83682250976562
$ws.Cells.Item(6, 5).Value = 0.99943125247955322
$ws.Cells.Item(7, 3).Value = 6004671
$ws.Cells.Item(7, 4).Value = 230.83682250976562
$ws.Cells.Item(7, 5).Value = 0.99990308284759521
$ws.Cells.Item(8, 3).Value = 6004671
$ws.Cells.Item(8, 4).Value = 230.83682250976562
$ws.Cells.Item(8, 5).Value = 0.9999840259552002
$ws.Cells.Item(9, 3).Value = 6004671
$ws.Cells.Item(9, 4).Value = 230.83682250976562
$ws.Cells.Item(9, 5).Value = 0.99999052286148071
$ws.Cells.Item(10, 3).Value = 6004671
$ws.Cells.Item(10, 4).Value = 230.83682250976562
$ws.Cells.Item(10, 5).Value = 0.99999964237213135

# Year 2005 (sheet index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 3).Value = 6005578
$ws.Cells.Item(2, 4).Value = 249.86701965332031
$ws.Cells.Item(2, 5).Value = 0.9662889838218689
$ws.Cells.Item(3, 3).Value = 6005578
$ws.Cells.Item(3, 4).Value = 249.86701965332031
$ws.Cells.Item(3, 5).Value = 0.970389723777771
$ws.Cells.Item(4, 3).Value = 6005578
$ws.Cells.Item(4, 4).Value = 249.86701965332031
$ws.Cells.Item(4, 5).Value = 0.97772204875946045
$ws.Cells.Item(5, 3).Value = 6005578
$ws.Cells.Item(5, 4).Value = 249.86701965332031
$ws.Cells.Item(5, 5).Value = 0.99678051471710205
$ws.Cells.Item(6, 3).Value = 6005578
$ws.Cells.Item(6, 4).Value = 249.86701965332031
$ws.Cells.Item(6, 5).Value = 0.99931961297988892
$ws.Cells.Item(7, 3).Value = 6005578
$ws.Cells.Item(7, 4).Value = 249.86701965332031
$ws.Cells.Item(7, 5).Value = 0.99987345933914185
$ws.Cells.Item(8, 3).Value = 6005578
$ws.Cells.Item(8, 4).Value = 249.86701965332031
$ws.Cells.Item(8, 5).Value = 0.99998182058334351
$ws.Cells.Item(9, 3).Value = 6005578
$ws.Cells.Item(9, 4).Value = 249.86701965332031
$ws.Cells.Item(9, 5).Value = 0.99998903274536133
$ws.Cells.Item(10, 3).Value = 6005578
$ws.Cells.Item(10, 4).Value = 249.86701965332031
$ws.Cells.Item(10, 5).Value = 0.99999934434890747
$ws.Cells.Item(11, 3).Value = 6005578
$ws.Cells.Item(11, 4).Value = 249.86701965332031
$ws.Cells.Item(11, 5).Value = 0.99999964237213135

# Year 2006 (sheet index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 3).Value = 6002319
$ws.Cells.Item(2, 4).Value = 315.19244384765625
$ws.Cells.Item(2, 5).Value = 0.96240085363388062
$ws.Cells.Item(3, 3).Value = 6002319
$ws.Cells.Item(3, 4).Value = 315.19244384765625
$ws.Cells.Item(3, 5).Value = 0.96693128347396851
$ws.Cells.Item(4, 3).Value = 6002319
$ws.Cells.Item(4, 4).Value = 315.19244384765625
$ws.Cells.Item(4, 5).Value = 0.97507262229919434
$ws.Cells.Item(5, 3).Value = 6002319
$ws.Cells.Item(5, 4).Value = 315.19244384765625
$ws.Cells.Item(5, 5).Value = 0.99639022350311279
$ws.Cells.Item(6, 3).Value = 6002319
$ws.Cells.Item(6, 4).Value = 315.19244384765625
$ws.Cells.Item(6, 5).Value = 0.99921178817749023
$ws.Cells.Item(7, 3).Value = 6002319
$ws.Cells.Item(7, 4).Value = 315.19244384765625
$ws.Cells.Item(7, 5).Value = 0.99984824657440186
$ws.Cells.Item(8, 3).Value = 6002319
$ws.Cells.Item(8, 4).Value = 315.19244384765625
$ws.Cells.Item(8, 5).Value = 0.99997466802597046
$ws.Cells.Item(9, 3).Value = 6002319
$ws.Cells.Item(9, 4).Value = 315.19244384765625
$ws.Cells.Item(9, 5).Value = 0.9999840259552002
$ws.Cells.Item(10, 3).Value = 6002319
$ws.Cells.Item(10, 4).Value = 315.19244384765625
$ws.Cells.Item(10, 5).Value = 0.99999785423278809
$ws.Cells.Item(11, 3).Value = 6002319
$ws.Cells.Item(11, 4).Value = 315.19244384765625
$ws.Cells.Item(11, 5).Value = 0.99999868869781494

# Year 2007 (sheet index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 3).Value = 6009824
$ws.Cells.Item(2, 4).Value = 323.35556030273438
$ws.Cells.Item(2, 5).Value = 0.95716977119445801
$ws.Cells.Item(3, 3).Value = 6009824
$ws.Cells.Item(3, 4).Value = 323.35556030273438
$ws.Cells.Item(3, 5).Value = 0.96256715059280396
$ws.Cells.Item(4, 3).Value = 6009824
$ws.Cells.Item(4, 4).Value = 323.35556030273438
$ws.Cells.Item(4, 5).Value = 0.9713473916053772
$ws.Cells.Item(5, 3).Value = 6009824
$ws.Cells.Item(5, 4).Value = 323.35556030273438
$ws.Cells.Item(5, 5).Value = 0.99594897031784058
$ws.Cells.Item(6, 3).Value = 6009824
$ws.Cells.Item(6, 4).Value = 323.35556030273438
$ws.Cells.Item(6, 5).Value = 0.99912875890731812
$ws.Cells.Item(7, 3).Value = 6009824
$ws.Cells.Item(7, 4).Value = 323.35556030273438
$ws.Cells.Item(7, 5).Value = 0.99983125925064087
$ws.Cells.Item(8, 3).Value = 6009824
$ws.Cells.Item(8, 4).Value = 323.35556030273438
$ws.Cells.Item(8, 5).Value = 0.9999728798866272
$ws.Cells.Item(9, 3).Value = 6009824
$ws.Cells.Item(9, 4).Value = 323.35556030273438
$ws.Cells.Item(9, 5).Value = 0.9999852180480957
$ws.Cells.Item(10, 3).Value = 6009824
$ws.Cells.Item(10, 4).Value = 323.35556030273438
$ws.Cells.Item(10, 5).Value = 0.99999868869781494
$ws.Cells.Item(11, 3).Value = 6009824
$ws.Cells.Item(11, 4).Value = 323.35556030273438
$ws.Cells.Item(11, 5).Value = 0.99999934434890747

# Year 2009 (sheet index 9)
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 3).Value = 6048279
$ws.Cells.Item(2, 4).Value = 360.6739501953125
$ws.Cells.Item(2, 5).Value = 0.95419061183929443
$ws.Cells.Item(3, 3).Value = 6048279
$ws.Cells.Item(3, 4).Value = 360.6739501953125
$ws.Cells.Item(3, 5).Value = 0.95903182029724121
$ws.Cells.Item(4, 3).Value = 6048279
$ws.Cells.Item(4, 4).Value = 360.6739501953125
$ws.Cells.Item(4, 5).Value = 0.96807044744491577
$ws.Cells.Item(5, 3).Value = 6048279
$ws.Cells.Item(5, 4).Value = 360.6739501953125
$ws.Cells.Item(5, 5).Value = 0.9952123761177063
$ws.Cells.Item(6, 3).Value = 6048279
$ws.Cells.Item(6, 4).Value = 360.6739501953125
$ws.Cells.Item(6, 5).Value = 0.99899905920028687
$ws.Cells.Item(7, 3).Value = 6048279
$ws.Cells.Item(7, 4).Value = 360.6739501953125
$ws.Cells.Item(7, 5).Value = 0.99980521202087402
$ws.Cells.Item(8, 3).Value = 6048279
$ws.Cells.Item(8, 4).Value = 360.6739501953125
$ws.Cells.Item(8, 5).Value = 0.99996626377105713
$ws.Cells.Item(9, 3).Value = 6048279
$ws.Cells.Item(9, 4).Value = 360.6739501953125
$ws.Cells.Item(9, 5).Value = 0.99997997283935547
$ws.Cells.Item(10, 3).Value = 6048279
$ws.Cells.Item(10, 4).Value = 360.6739501953125

# Year 2010 (sheet index 10)
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(2, 3).Value = 6068249
$ws.Cells.Item(2, 4).Value = 389.43756103515625
$ws.Cells.Item(2, 5).Value = 0.95038139820098877
$ws.Cells.Item(3, 3).Value = 6068249
$ws.Cells.Item(3, 4).Value = 389.43756103515625
$ws.Cells.Item(3, 5).Value = 0.95579665899276733
$ws.Cells.Item(4, 3).Value = 6068249
$ws.Cells.Item(4, 4).Value = 389.43756103515625
$ws.Cells.Item(4, 5).Value = 0.96584975719451904
$ws.Cells.Item(5, 3).Value = 6068249
$ws.Cells.Item(5, 4).Value = 389.43756103515625
$ws.Cells.Item(5, 5).Value = 0.99484330415725708
$ws.Cells.Item(6, 3).Value = 6068249
$ws.Cells.Item(6, 4).Value = 389.43756103515625
$ws.Cells.Item(6, 5).Value = 0.99888551235198975
$ws.Cells.Item(7, 3).Value = 6068249
$ws.Cells.Item(7, 4).Value = 389.43756103515625
$ws.Cells.Item(7, 5).Value = 0.99978101253509521
$ws.Cells.Item(8, 3).Value = 6068249
$ws.Cells.Item(8, 4).Value = 389.43756103515625
$ws.Cells.Item(8, 5).Value = 0.99995797872543335
$ws.Cells.Item(9, 3).Value = 6068249
$ws.Cells.Item(9, 4).Value = 389.43756103515625
$ws.Cells.Item(9, 5).Value = 0.99997591972351074
$ws.Cells.Item(10, 3).Value = 6068249
$ws.Cells.Item(10, 4).Value = 389.43756103515625
$ws.Cells.Item(11, 3).Value = 6068249
$ws.Cells.Item(11, 4).Value = 389.43756103515625

# Year 2012 (sheet index 11)
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(2, 3).Value = 6113975
$ws.Cells.Item(2, 4).Value = 239.74188232421875
$ws.Cells.Item(2, 5).Value = 0.97683942317962646
$ws.Cells.Item(3, 3).Value = 6113975
$ws.Cells.Item(3, 4).Value = 239.74188232421875
$ws.Cells.Item(3, 5).Value = 0.97803264856338501
$ws.Cells.Item(4, 3).Value = 6113975
$ws.Cells.Item(4, 4).Value = 239.74188232421875
$ws.Cells.Item(4, 5).Value = 0.98100060224533081
$ws.Cells.Item(5, 3).Value = 6113975
$ws.Cells.Item(5, 4).Value = 239.74188232421875
$ws.Cells.Item(5, 5).Value = 0.99575364589691162
$ws.Cells.Item(6, 3).Value = 6113975
$ws.Cells.Item(6, 4).Value = 239.74188232421875
$ws.Cells.Item(6, 5).Value = 0.99887293577194214
$ws.Cells.Item(7, 3).Value = 6113975
$ws.Cells.Item(7, 4).Value = 239.74188232421875
$ws.Cells.Item(7, 5).Value = 0.9997219443321228
$ws.Cells.Item(8, 3).Value = 6113975
$ws.Cells.Item(8, 4).Value = 239.74188232421875
$ws.Cells.Item(8, 5).Value = 0.9999430775642395
$ws.Cells.Item(9, 3).Value = 6113975
$ws.Cells.Item(9, 4).Value = 239.74188232421875
$ws.Cells.Item(9, 5).Value = 0.99997025728225708
$ws.Cells.Item(10, 3).Value = 6113975
$ws.Cells.Item(10, 4).Value = 239.74188232421875

# Year 2013 (sheet index 12)
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 3).Value = 6138839
$ws.Cells.Item(2, 4).Value = 280.01898193359375
$ws.Cells.Item(2, 5).Value = 0.97289276123046875
$ws.Cells.Item(3, 3).Value = 6138839
$ws.Cells.Item(3, 4).Value = 280.01898193359375
$ws.Cells.Item(3, 5).Value = 0.97424709796905518
$ws.Cells.Item(4, 3).Value = 6138839
$ws.Cells.Item(4, 4).Value = 280.01898193359375
$ws.Cells.Item(4, 5).Value = 0.97784191370010376
$ws.Cells.Item(5, 3).Value = 6138839
$ws.Cells.Item(5, 4).Value = 280.01898193359375
$ws.Cells.Item(5, 5).Value = 0.9950181245803833
$ws.Cells.Item(6, 3).Value = 6138839
$ws.Cells.Item(6, 4).Value = 280.01898193359375
$ws.Cells.Item(6, 5).Value = 0.99870824813842773
$ws.Cells.Item(7, 3).Value = 6138839
$ws.Cells.Item(7, 4).Value = 280.01898193359375
$ws.Cells.Item(7, 5).Value = 0.99969702959060669
$ws.Cells.Item(8, 3).Value = 6138839
$ws.Cells.Item(8, 4).Value = 280.01898193359375
$ws.Cells.Item(8, 5).Value = 0.99993419647216797
$ws.Cells.Item(9, 3).Value = 6138839
$ws.Cells.Item(9, 4).Value = 280.01898193359375
$ws.Cells.Item(9, 5).Value = 0.99996399879455566
$ws.Cells.Item(10, 3).Value = 6138839
$ws.Cells.Item(10, 4).Value = 280.01898193359375

# Year 2014 (sheet index 13)
$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(2, 3).Value = 6162955
$ws.Cells.Item(2, 4).Value = 325.3929443359375
$ws.Cells.Item(2, 5).Value = 0.9684099555015564
$ws.Cells.Item(3, 3).Value = 6162955
$ws.Cells.Item(3, 4).Value = 325.3929443359375
$ws.Cells.Item(3, 5).Value = 0.96977150440216064
$ws.Cells.Item(4, 3).Value = 6162955
$ws.Cells.Item(4, 4).Value = 325.3929443359375
$ws.Cells.Item(4, 5).Value = 0.97381728887557983
$ws.Cells.Item(5, 3).Value = 6162955
$ws.Cells.Item(5, 4).Value = 325.3929443359375
$ws.Cells.Item(5, 5).Value = 0.99428242444992065
$ws.Cells.Item(6, 3).Value = 6162955
$ws.Cells.Item(6, 4).Value = 325.3929443359375
$ws.Cells.Item(6, 5).Value = 0.99857747554779053
$ws.Cells.Item(7, 3).Value = 6162955
$ws.Cells.Item(7, 4).Value = 325.3929443359375
$ws.Cells.Item(7, 5).Value = 0.999672532081604
$ws.Cells.Item(8, 3).Value = 6162955
$ws.Cells.Item(8, 4).Value = 325.3929443359375
$ws.Cells.Item(8, 5).Value = 0.9999312162399292
$ws.Cells.Item(9, 3).Value = 6162955
$ws.Cells.Item(9, 4).Value = 325.3929443359375
$ws.Cells.Item(9, 5).Value = 0.99995958805084229
$ws.Cells.Item(10, 3).Value = 6162955
$ws.Cells.Item(10, 4).Value = 325.3929443359375
$ws.Cells.Item(10, 5).Value = 0.99999850988388062

# Year 2015 (sheet index 14)
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(2, 3).Value = 6183676
$ws.Cells.Item(2, 4).Value = 393.77197265625
$ws.Cells.Item(2, 5).Value = 0.96501660346984863
$ws.Cells.Item(3, 3).Value = 6183676
$ws.Cells.Item(3, 4).Value = 393.77197265625
$ws.Cells.Item(3, 5).Value = 0.96644598245620728
$ws.Cells.Item(4, 3).Value = 6183676
$ws.Cells.Item(4, 4).Value = 393.77197265625
$ws.Cells.Item(4, 5).Value = 0.97067344188690186
$ws.Cells.Item(5, 3).Value = 6183676
$ws.Cells.Item(5, 4).Value = 393.77197265625
$ws.Cells.Item(5, 5).Value = 0.99345391988754272
$ws.Cells.Item(6, 3).Value = 6183676
$ws.Cells.Item(6, 4).Value = 393.77197265625
$ws.Cells.Item(6, 5).Value = 0.99839949607849121
$ws.Cells.Item(7, 3).Value = 6183676
$ws.Cells.Item(7, 4).Value = 393.77197265625
$ws.Cells.Item(7, 5).Value = 0.99965053796768188
$ws.Cells.Item(8, 3).Value = 6183676
$ws.Cells.Item(8, 4).Value = 393.77197265625
$ws.Cells.Item(8, 5).Value = 0.99992883205413818
$ws.Cells.Item(9, 3).Value = 6183676
$ws.Cells.Item(9, 4).Value = 393.77197265625
$ws.Cells.Item(9, 5).Value = 0.99995732307434082
$ws.Cells.Item(10, 3).Value = 6183676
$ws.Cells.Item(10, 4).Value = 393.77197265625
$ws.Cells.Item(10, 5).Value = 0.99999791383743286
$ws.Cells.Item(11, 3).Value = 6183676
$ws.Cells.Item(11, 4).Value = 393.77197265625

# Year 2016 (sheet index 15)
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(2, 3).Value = 6200800
$ws.Cells.Item(2, 4).Value = 368.2369384765625
$ws.Cells.Item(2, 5).Value = 0.96577620506286621
$ws.Cells.Item(3, 3).Value = 6200800
$ws.Cells.Item(3, 4).Value = 368.2369384765625
$ws.Cells.Item(3, 5).Value = 0.96695071458816528
$ws.Cells.Item(4, 3).Value = 6200800
$ws.Cells.Item(4, 4).Value = 368.2369384765625
$ws.Cells.Item(4, 5).Value = 0.97109532356262207
$ws.Cells.Item(5, 3).Value = 6200800
$ws.Cells.Item(5, 4).Value = 368.2369384765625
$ws.Cells.Item(5, 5).Value = 0.99311429262161255
$ws.Cells.Item(6, 3).Value = 6200800
$ws.Cells.Item(6, 4).Value = 368.2369384765625
$ws.Cells.Item(6, 5).Value = 0.99829262495040894
$ws.Cells.Item(7, 3).Value = 6200800
$ws.Cells.Item(7, 4).Value = 368.2369384765625
$ws.Cells.Item(7, 5).Value = 0.9996299147605896
$ws.Cells.Item(8, 3).Value = 6200800
$ws.Cells.Item(8, 4).Value = 368.2369384765625
$ws.Cells.Item(8, 5).Value = 0.99992680549621582
$ws.Cells.Item(9, 3).Value = 6200800
$ws.Cells.Item(9, 4).Value = 368.2369384765625
$ws.Cells.Item(9, 5).Value = 0.99995642900466919
$ws.Cells.Item(10, 3).Value = 6200800
$ws.Cells.Item(10, 4).Value = 368.2369384765625
$ws.Cells.Item(10, 5).Value = 0.99999630451202393
$ws.Cells.Item(11, 3).Value = 6200800
$ws.Cells.Item(11, 4).Value = 368.2369384765625
$ws.Cells.Item(11, 5).Value = 0.99999886751174927

# Year 2017 (sheet index 16)
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(2, 3).Value = 6213533
$ws.Cells.Item(2, 4).Value = 400.72366333007812
$ws.Cells.Item(2, 5).Value = 0.96211367845535278
$ws.Cells.Item(3, 3).Value = 6213533
$ws.Cells.Item(3, 4).Value = 400.72366333007812
$ws.Cells.Item(3, 5).Value = 0.96327805519104004
$ws.Cells.Item(4, 3).Value = 6213533
$ws.Cells.Item(4, 4).Value = 400.72366333007812
$ws.Cells.Item(4, 5).Value = 0.96767377853393555
$ws.Cells.Item(5, 3).Value = 6213533
$ws.Cells.Item(5, 4).Value = 400.72366333007812
$ws.Cells.Item(5, 5).Value = 0.9926491379737854
$ws.Cells.Item(6, 3).Value = 6213533
$ws.Cells.Item(6, 4).Value = 400.72366333007812
$ws.Cells.Item(6, 5).Value = 0.99822080135345459
$ws.Cells.Item(7, 3).Value = 6213533
$ws.Cells.Item(7, 4).Value = 400.72366333007812
$ws.Cells.Item(7, 5).Value = 0.99962180852890015
$ws.Cells.Item(8, 3).Value = 6213533
$ws.Cells.Item(8, 4).Value = 400.72366333007812
$ws.Cells.Item(8, 5).Value = 0.99992567300796509
$ws.Cells.Item(9, 3).Value = 6213533
$ws.Cells.Item(9, 4).Value = 400.72366333007812
$ws.Cells.Item(9, 5).Value = 0.99995815753936768
$ws.Cells.Item(10, 3).Value = 6213533
$ws.Cells.Item(10, 4).Value = 400.72366333007812
$ws.Cells.Item(10, 5).Value = 0.99999821186065674
$ws.Cells.Item(11, 3).Value = 6213533
$ws.Cells.Item(11, 4).Value = 400.72366333007812
